# Apply the "added has posted and legal_entites" edit to the vocab-tagged
# workbook.
#
# Summary of the change:
#  - Column D header changes from "assets" to "collateral"
#  - A new column E is added with header "legal_entity" and values
#    cgmi / cgml / cgma for rows 2-4
#  - A new value "has posted" is added to column A, row 4 (it already
#    existed in column C, row 4)
#  - The active selection moves from A4 to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D header: assets -> collateral
$ws.Range("D1").Value = "collateral"

# New column E: legal_entity header + values
$ws.Range("E1").Value = "legal_entity"
$ws.Range("E2").Value = "cgmi"
$ws.Range("E3").Value = "cgml"
$ws.Range("E4").Value = "cgma"

# Column A gains "has posted" in row 4 (already present in C4)
$ws.Range("A4").Value = "has posted"

# Update the selected cell to match the saved workbook state
$ws.Range("E5").Select()
